$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$notes = $s.NotesPage
$notesShape = $notes.Shapes.AddPlaceholder(2)
$notesShape.TextFrame.TextRange.Text = "This video will discuss administrative simplification, a key goal of the Massachusetts All-Payer Claims Database."
